$wb = $excel.ActiveWorkbook
$survey = $wb.Worksheets.Item(1)
$choices = $wb.Worksheets.Item(2)

# --- choices sheet: add "photo_review" and "whatsapp" options to the
#     followup_request choice list, right before the "none" option ---
$choices.Rows("35:36").Insert()
$choices.Range("A35").Value = "followup_request"
$choices.Range("B35").Value = "photo_review"
$choices.Range("C35").Value = "Photo Review"
$choices.Range("A36").Value = "followup_request"
$choices.Range("B36").Value = "whatsapp"
$choices.Range("C36").Value = "Whatsapp"

# --- survey sheet: row 22 "followup_request" becomes a multi-select question ---
$survey.Range("A22").Value = "select_multiple followup_request"

# --- restore view/selection state on the touched sheets ---
$choices.Activate()
$choices.Range("E40").Select()

$survey.Activate()
$survey.Range("C18").Select()
